$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 81.05837566666666
$ws.Range("H2").Value = 243.175127
$ws.Range("I2").Value = 0.3545816884225585
$ws.Range("J2").Value = 0.3545816884225585
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 15.35884066666667
$ws.Range("N2").Value = 46.076522
$ws.Range("O2").Value = 0.1012042817263867
$ws.Range("P2").Value = 0.1012042817263867
$ws.Range("Q2").Value = 1244.962676563144
$ws.Range("R2").Value = 11204.66408906829
$ws.Range("S2").Value = 0.03588518509013446
$ws.Range("T2").Value = 0.03588518509013447

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 81.05837566666666
$ws.Range("H3").Value = 243.175127
$ws.Range("I3").Value = 0.3545816884225585
$ws.Range("J3").Value = 0.3545816884225585
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 50.59256466666667
$ws.Range("N3").Value = 151.777694
$ws.Range("O3").Value = 0.3333704853712116
$ws.Range("P3").Value = 0.3333704853712116
$ws.Range("Q3").Value = 4100.951112690793
$ws.Range("R3").Value = 36908.56001421713
$ws.Range("S3").Value = 0.118207069573172
$ws.Range("T3").Value = 0.118207069573172

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 81.05837566666666
$ws.Range("H4").Value = 243.175127
$ws.Range("I4").Value = 0.3545816884225585
$ws.Range("J4").Value = 0.3545816884225585
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 60.37715666666667
$ws.Range("N4").Value = 181.13147
$ws.Range("O4").Value = 0.397844271305776
$ws.Range("P4").Value = 0.397844271305776
$ws.Range("Q4").Value = 4894.074246771855
$ws.Range("R4").Value = 44046.66822094669
$ws.Range("S4").Value = 0.1410682934488445
$ws.Range("T4").Value = 0.1410682934488445

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 81.05837566666666
$ws.Range("H5").Value = 243.175127
$ws.Range("I5").Value = 0.3545816884225585
$ws.Range("J5").Value = 0.3545816884225585
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 25.43221733333333
$ws.Range("N5").Value = 76.29665199999999
$ws.Range("O5").Value = 0.1675809615966257
$ws.Range("P5").Value = 0.1675809615966258
$ws.Range("Q5").Value = 2061.494226641645
$ws.Range("R5").Value = 18553.4480397748
$ws.Range("S5").Value = 0.05942114031040749
$ws.Range("T5").Value = 0.0594211403104075

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 123.018252
$ws.Range("H6").Value = 369.054756
$ws.Range("I6").Value = 0.5381309351710768
$ws.Range("J6").Value = 0.5381309351710768
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 15.35884066666667
$ws.Range("N6").Value = 46.076522
$ws.Range("O6").Value = 0.1012042817263867
$ws.Range("P6").Value = 0.1012042817263867
$ws.Range("Q6").Value = 1889.417731559848
$ws.Range("R6").Value = 17004.75958403863
$ws.Range("S6").Value = 0.05446115476873756
$ws.Range("T6").Value = 0.05446115476873758

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 123.018252
$ws.Range("H7").Value = 369.054756
$ws.Range("I7").Value = 0.5381309351710768
$ws.Range("J7").Value = 0.5381309351710768
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 50.59256466666667
$ws.Range("N7").Value = 151.777694
$ws.Range("O7").Value = 0.3333704853712116
$ws.Range("P7").Value = 0.3333704853712116
$ws.Range("Q7").Value = 6223.808869490296
$ws.Range("R7").Value = 56014.27982541266
$ws.Range("S7").Value = 0.1793969710512459
$ws.Range("T7").Value = 0.1793969710512459

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 123.018252
$ws.Range("H8").Value = 369.054756
$ws.Range("I8").Value = 0.5381309351710768
$ws.Range("J8").Value = 0.5381309351710768
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 60.37715666666667
$ws.Range("N8").Value = 181.13147
$ws.Range("O8").Value = 0.397844271305776
$ws.Range("P8").Value = 0.397844271305776
$ws.Range("Q8").Value = 7427.49227386348
$ws.Range("R8").Value = 66847.43046477133
$ws.Range("S8").Value = 0.2140923097702328
$ws.Range("T8").Value = 0.2140923097702328

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 123.018252
$ws.Range("H9").Value = 369.054756
$ws.Range("I9").Value = 0.5381309351710768
$ws.Range("J9").Value = 0.5381309351710768
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 25.43221733333333
$ws.Range("N9").Value = 76.29665199999999
$ws.Range("O9").Value = 0.1675809615966257
$ws.Range("P9").Value = 0.1675809615966258
$ws.Range("Q9").Value = 3128.626920830768
$ws.Range("R9").Value = 28157.64228747691
$ws.Range("S9").Value = 0.09018049958086052
$ws.Range("T9").Value = 0.09018049958086054

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.3624666666666667
$ws.Range("H10").Value = 1.0874
$ws.Range("I10").Value = 0.001585573873230423
$ws.Range("J10").Value = 0.001585573873230423
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 15.35884066666667
$ws.Range("N10").Value = 46.076522
$ws.Range("O10").Value = 0.1012042817263867
$ws.Range("P10").Value = 0.1012042817263867
$ws.Range("Q10").Value = 5.567067780311111
$ws.Range("R10").Value = 50.10361002279999
$ws.Range("S10").Value = 0.0001604668649644098
$ws.Range("T10").Value = 0.0001604668649644099

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.3624666666666667
$ws.Range("H11").Value = 1.0874
$ws.Range("I11").Value = 0.001585573873230423
$ws.Range("J11").Value = 0.001585573873230423
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 50.59256466666667
$ws.Range("N11").Value = 151.777694
$ws.Range("O11").Value = 0.3333704853712116
$ws.Range("P11").Value = 0.3333704853712116
$ws.Range("Q11").Value = 18.33811827284444
$ws.Range("R11").Value = 165.0430644556
$ws.Range("S11").Value = 0.0005285835317107381
$ws.Range("T11").Value = 0.0005285835317107381

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.3624666666666667
$ws.Range("H12").Value = 1.0874
$ws.Range("I12").Value = 0.001585573873230423
$ws.Range("J12").Value = 0.001585573873230423
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 60.37715666666667
$ws.Range("N12").Value = 181.13147
$ws.Range("O12").Value = 0.397844271305776
$ws.Range("P12").Value = 0.397844271305776
$ws.Range("Q12").Value = 21.88470671977778
$ws.Range("R12").Value = 196.962360478
$ws.Range("S12").Value = 0.0006308114821968347
$ws.Range("T12").Value = 0.0006308114821968347

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.3624666666666667
$ws.Range("H13").Value = 1.0874
$ws.Range("I13").Value = 0.001585573873230423
$ws.Range("J13").Value = 0.001585573873230423
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 25.43221733333333
$ws.Range("N13").Value = 76.29665199999999
$ws.Range("O13").Value = 0.1675809615966257
$ws.Range("P13").Value = 0.1675809615966258
$ws.Range("Q13").Value = 9.218331042755555
$ws.Range("R13").Value = 82.96497938479999
$ws.Range("S13").Value = 0.0002657119943584407
$ws.Range("T13").Value = 0.0002657119943584408

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 24.16373066666667
$ws.Range("H14").Value = 72.491192
$ws.Range("I14").Value = 0.1057018025331343
$ws.Range("J14").Value = 0.1057018025331344
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 15.35884066666667
$ws.Range("N14").Value = 46.076522
$ws.Range("O14").Value = 0.1012042817263867
$ws.Range("P14").Value = 0.1012042817263867
$ws.Range("Q14").Value = 371.1268892215804
$ws.Range("R14").Value = 3340.142002994224
$ws.Range("S14").Value = 0.01069747500255022
$ws.Range("T14").Value = 0.01069747500255022

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 24.16373066666667
$ws.Range("H15").Value = 72.491192
$ws.Range("I15").Value = 0.1057018025331343
$ws.Range("J15").Value = 0.1057018025331344
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 50.59256466666667
$ws.Range("N15").Value = 151.777694
$ws.Range("O15").Value = 0.3333704853712116
$ws.Range("P15").Value = 0.3333704853712116
$ws.Range("Q15").Value = 1222.50510634125
$ws.Range("R15").Value = 11002.54595707125
$ws.Range("S15").Value = 0.03523786121508295
$ws.Range("T15").Value = 0.03523786121508296

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 24.16373066666667
$ws.Range("H16").Value = 72.491192
$ws.Range("I16").Value = 0.1057018025331343
$ws.Range("J16").Value = 0.1057018025331344
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 60.37715666666667
$ws.Range("N16").Value = 181.13147
$ws.Range("O16").Value = 0.397844271305776
$ws.Range("P16").Value = 0.397844271305776
$ws.Range("Q16").Value = 1458.937352112471
$ws.Range("R16").Value = 13130.43616901224
$ws.Range("S16").Value = 0.04205285660450186
$ws.Range("T16").Value = 0.04205285660450186

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 24.16373066666667
$ws.Range("H17").Value = 72.491192
$ws.Range("I17").Value = 0.1057018025331343
$ws.Range("J17").Value = 0.1057018025331344
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 25.43221733333333
$ws.Range("N17").Value = 76.29665199999999
$ws.Range("O17").Value = 0.1675809615966257
$ws.Range("P17").Value = 0.1675809615966258
$ws.Range("Q17").Value = 614.5372498987981
$ws.Range("R17").Value = 5530.835249089184
$ws.Range("S17").Value = 0.0177136097109993
$ws.Range("T17").Value = 0.01771360971099931
